$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("An introspective look at the relationship between Hawking and the space/time contingent. This film explores the Gallilean and Newtonian laws and ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "An introspective look at the relationship between Hawking and the space/time contingent. ## This film explores the Gallilean and Newtonian laws and " } else { Write-Output "NOT FOUND: 0" }
$rng = $d.Content
$found = $rng.Find.Execute(" relation to Einstein's Theory of General Relativity. The film is methodically directed, exposing details of the man (Hawking) as well as his work (Black Holes). Interviews with his family are a little too long so sadly there is less development of his theories and ideas.  A Philip Glass soundtrack superbly compliments the film. Only one other man could compose such haunting instellar melodies (Jean Michel Jarre). Overall I would highly recommend this movie on the basis of Hawking's 'nuggets of wisdom' and his adequate explanation of an Event Horizon!", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " relation to Einstein's Theory of General Relativity. ## The film is methodically directed, exposing details of the man (Hawking) as well as his work (Black Holes). ## Interviews with his family are a little too long so sadly there is less development of his theories and ideas. ## A Philip Glass soundtrack superbly compliments the film. ## Only one other man could compose such haunting instellar melodies (Jean Michel Jarre). ## Overall I would highly recommend this movie on the basis of Hawking's 'nuggets of wisdom' and his adequate explanation of an Event Horizon!" } else { Write-Output "NOT FOUND: 1" }
$rng = $d.Content
$found = $rng.Find.Execute("/vremenskog kontingenta. Ovaj film ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "/vremenskog kontingenta. ## Ovaj film " } else { Write-Output "NOT FOUND: 2" }
$rng = $d.Content
$found = $rng.Find.Execute(" odnos prema Einsteinovoj teoriji opće relativnosti. Film je metodički režiran, otkrivajući detalje ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " odnos prema Einsteinovoj teoriji opće relativnosti. ## Film je metodički režiran, otkrivajući detalje " } else { Write-Output "NOT FOUND: 3" }
$rng = $d.Content
$found = $rng.Find.Execute(" (Hawking) kao i njegov rad (Crne rupe). Intervjui sa", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " (Hawking) kao i njegov rad (Crne rupe). ## Intervjui sa" } else { Write-Output "NOT FOUND: 4" }
$rng = $d.Content
$found = $rng.Find.Execute(". Philip Glass ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = ". ## Philip Glass " } else { Write-Output "NOT FOUND: 5" }
$rng = $d.Content
$found = $rng.Find.Execute(" filmu. Samo je jedan drugi čovjek mogao skladati takve uklete instelarne melodije (Jean Michel Jarre). Sveukupno bih ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " filmu. ## Samo je jedan drugi čovjek mogao skladati takve uklete instelarne melodije (Jean Michel Jarre). ## Sveukupno bih " } else { Write-Output "NOT FOUND: 6" }
$rng = $d.Content
$found = $rng.Find.Execute("My god...i have not seen such an awful movie in a long...long time...saw it last night and wanted to leave after 20 minutes...keira knightley tries really really hard in this one, but she cant handle it..dropped her accent every once in a while and didn't have the charisma to fill the role...sienna millers acting gets you to a point where you start to ask yourself: Has she ever had acting lessons? judging by the edge of love shes never been to acting class, but should consider to go in the near future...they both look really pretty..maybe thats what they should focus on in their future career..if they can be actresses everybody can!", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "My god...i have not seen such an awful movie in a long...long time… ## saw it last night and wanted to leave after 20 minutes… ## keira knightley tries really really hard in this one, but she cant handle it..dropped her accent every once in a while and didn't have the charisma to fill the role… ## sienna millers acting gets you to a point where you start to ask yourself: Has she ever had acting lessons? ## judging by the edge of love shes never been to acting class, but should consider to go in the near future…## they both look really pretty..maybe thats what they should focus on in their future career.. ## if they can be actresses everybody can!" } else { Write-Output "NOT FOUND: 7" }
$rng = $d.Content
$found = $rng.Find.Execute("Moj Bože... nisam vidio tako grozan film u dugo... dugo vremena... ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "Moj Bože... nisam vidio tako grozan film u dugo... dugo vremena... ## " } else { Write-Output "NOT FOUND: 8" }
$rng = $d.Content
$found = $rng.Find.Execute("sam ga sinoć i htio otići nakon 20 minuta... Keira Knightley pokušava stvarno jako ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "sam ga sinoć i htio otići nakon 20 minuta... ## Keira Knightley pokušava stvarno jako " } else { Write-Output "NOT FOUND: 9" }
$rng = $d.Content
$found = $rng.Find.Execute("karizmu ispuniti ulogu... Sienna Millers ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "karizmu ispuniti ulogu... ## Sienna Millers " } else { Write-Output "NOT FOUND: 10" }
$rng = $d.Content
$found = $rng.Find.Execute(": Je li ikad imala satove glume? Sudeći po ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = ": Je li ikad imala satove glume? ## Sudeći po " } else { Write-Output "NOT FOUND: 11" }
$rng = $d.Content
$found = $rng.Find.Execute("u bliskoj budućnosti... ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "u bliskoj budućnosti… ## " } else { Write-Output "NOT FOUND: 12" }
$rng = $d.Content
$found = $rng.Find.Execute(".. ako oni mogu biti glumice svatko može!", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = ".. ## ako oni mogu biti glumice svatko može!" } else { Write-Output "NOT FOUND: 13" }
$rng = $d.Content
$found = $rng.Find.Execute("Thanks to other reviewers who directed me to this product when I was told I was anemic. Now been taking these for about 4 months and the anemia is gone.  Good product.  Easily digested (unlike some other iron supplements).", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "Thanks to other reviewers who directed me to this product when I was told I was anemic. ## Now been taking these for about 4 months and the anemia is gone. ## Good product. ## Easily digested (unlike some other iron supplements)." } else { Write-Output "NOT FOUND: 14" }
$rng = $d.Content
$found = $rng.Find.Execute("Zahvaljujući drugim recenzentima koji su me usmjerili na ovaj proizvod kad mi je rečeno da sam anemična. Sada uzimam ove ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "Zahvaljujući drugim recenzentima koji su me usmjerili na ovaj proizvod kad mi je rečeno da sam anemična. ## Sada uzimam ove " } else { Write-Output "NOT FOUND: 15" }
$rng = $d.Content
$found = $rng.Find.Execute(" oko 4 mjeseca, a anemija je nestala. Dobar proizvod. Jednostavno ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " oko 4 mjeseca, a anemija je nestala. ## Dobar proizvod. ## Jednostavno " } else { Write-Output "NOT FOUND: 16" }
$rng = $d.Content
$found = $rng.Find.Execute("This is one of my favorite desserts, and melts quickly in the mouth. This brand is good and it shipped well-packaged. Everyone should try this once. The amazon price is much better than the ones you find at science fairs.", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "This is one of my favorite desserts, and melts quickly in the mouth. ## This brand is good and it shipped well-packaged. ## Everyone should try this once. ## The amazon price is much better than the ones you find at science fairs." } else { Write-Output "NOT FOUND: 17" }
$rng = $d.Content
$found = $rng.Find.Execute("i brzo se topi u ustima. Ova marka je dobra i ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "i brzo se topi u ustima. ## Ova marka je dobra i " } else { Write-Output "NOT FOUND: 18" }
$rng = $d.Content
$found = $rng.Find.Execute("dobro zapakirana. Svatko bi trebao probati ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "dobro zapakirana. ## Svatko bi trebao probati " } else { Write-Output "NOT FOUND: 19" }
$rng = $d.Content
$found = $rng.Find.Execute("jednom. Cijena ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "jednom. ## Cijena " } else { Write-Output "NOT FOUND: 20" }
$rng = $d.Content
$found = $rng.Find.Execute("This is a fantastic puzzle/gift for young AND old. It is 32 triangular strong magnetic pieces that can fit together in a wide number of ways.  It's just great and you'll have trouble keeping it away from the adults.", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "This is a fantastic puzzle/gift for young AND old. ## It is 32 triangular strong magnetic pieces that can fit together in a wide number of ways. ## It's just great and you'll have trouble keeping it away from the adults." } else { Write-Output "NOT FOUND: 21" }
$rng = $d.Content
$found = $rng.Find.Execute("Ovo je fantastična zagonetka/poklon za mlade i stare. ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "Ovo je fantastična zagonetka/poklon za mlade i stare. ## " } else { Write-Output "NOT FOUND: 22" }
$rng = $d.Content
$found = $rng.Find.Execute(" broj načina. Sjajno je i teško ćeš ga držati podalje od odraslih.", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " broj načina. ## Sjajno je i teško ćeš ga držati podalje od odraslih." } else { Write-Output "NOT FOUND: 23" }
$rng = $d.Content
$found = $rng.Find.Execute("It's another bad zombie movie. Compared to the majority of thhe others, the only difference here is the main character is a female. The plot is the same. The action scenes are not engaging. Special effects  are so so.", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "It's another bad zombie movie. ## Compared to the majority of thhe others, the only difference here is the main character is a female. ## The plot is the same. ## The action scenes are not engaging. ## Special effects  are so so." } else { Write-Output "NOT FOUND: 24" }
$rng = $d.Content
$found = $rng.Find.Execute("Još jedan loš zombi film. U usporedbi s većinom drugih, jedina razlika ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "Još jedan loš zombi film. ## U usporedbi s većinom drugih, jedina razlika " } else { Write-Output "NOT FOUND: 25" }
$rng = $d.Content
$found = $rng.Find.Execute("žensko. ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "žensko. ## " } else { Write-Output "NOT FOUND: 26" }
$rng = $d.Content
$found = $rng.Find.Execute("je ista. Akcijske scene nisu zanimljive. Specijalni efekti su ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "je ista. ## Akcijske scene nisu zanimljive. ## Specijalni efekti su " } else { Write-Output "NOT FOUND: 27" }
$rng = $d.Content
$found = $rng.Find.Execute("YoYo seems out of balance. No matter ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "YoYo seems out of balance. ## No matter " } else { Write-Output "NOT FOUND: 28" }
$rng = $d.Content
$found = $rng.Find.Execute(" tilt to one side. Made it difficult", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " tilt to one side. ## Made it difficult" } else { Write-Output "NOT FOUND: 29" }
$rng = $d.Content
$found = $rng.Find.Execute(" or do tricks. I have a little ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " or do tricks.  ## I have a little " } else { Write-Output "NOT FOUND: 30" }
$rng = $d.Content
$found = $rng.Find.Execute("Yoyo izgleda izvan ravnoteže. Bez obzira ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "Yoyo izgleda izvan ravnoteže. ## Bez obzira " } else { Write-Output "NOT FOUND: 31" }
$rng = $d.Content
$found = $rng.Find.Execute(" na jednu stranu. Otežalo je", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " na jednu stranu. ## Otežalo je" } else { Write-Output "NOT FOUND: 32" }
$rng = $d.Content
$found = $rng.Find.Execute(" ili raditi trikove. Imam malo iskustva s ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = " ili raditi trikove. ## Imam malo iskustva s " } else { Write-Output "NOT FOUND: 33" }
$rng = $d.Content
$found = $rng.Find.Execute("If you wrap your wrists properly, you'll see these are both too narrow and too short, way too short. Do not get these if you are hitting the heavy bag. They just won't protect/support your wrists or knuckles.", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "If you wrap your wrists properly, you'll see these are both too narrow and too short, way too short. ## Do not get these if you are hitting the heavy bag. ## They just won't protect/support your wrists or knuckles." } else { Write-Output "NOT FOUND: 34" }
$rng = $d.Content
$found = $rng.Find.Execute("prekratka. Nemojte ih ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = "prekratka. ## Nemojte ih " } else { Write-Output "NOT FOUND: 35" }
$rng = $d.Content
$found = $rng.Find.Execute(". Oni jednostavno neće štititi/podupirati vaše zglobove ili ", $true, $false, $false, $false, $false, $true, 1, $false)
if ($found) { $rng.Text = ". ## Oni jednostavno neće štititi/podupirati vaše zglobove ili " } else { Write-Output "NOT FOUND: 36" }
